$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'313.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.47%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'38.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.14%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.075"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.44%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07769"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-4.85%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.353"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.64%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.907"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-4.21%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'8.187"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9178"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.94%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1243"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-4.51%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'-3.56%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08830"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.55%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03388"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-3.88%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09702"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.34%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001377"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.14%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005871"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-9.86%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.533"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.71%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.981"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-4.63%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-1.77%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1296"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.47%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.029"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.78%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'4.09%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'5,589.85%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04399"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.94%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-2.26%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-10.79%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001350"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-65.32%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02135"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-4.71%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04975"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-4.22%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007744"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.04%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009882"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-4.06%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-3.75%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-1.97%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009676"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'9.25%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006515"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.47%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003074"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'2.20%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.12%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.01%"
$ws.Range("E51").Style = "Normal"
